$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 37
$ws.Range("F3").Value = 821
$ws.Range("F4").Value = 2468
$ws.Range("F5").Value = 524
$ws.Range("F6").Value = 447
$ws.Range("F7").Value = 269
$ws.Range("F10").Value = 1147
$ws.Range("F11").Value = 507
$ws.Range("F12").Value = 270
$ws.Range("F15").Value = 5250
$ws.Range("F17").Value = 1503
$ws.Range("F18").Value = 3851
$ws.Range("F22").Value = 4286
$ws.Range("F23").Value = 5680
$ws.Range("F24").Value = 136
$ws.Range("F25").Value = 1001
$ws.Range("F26").Value = 615
$ws.Range("F27").Value = 3554
$ws.Range("F28").Value = 436
$ws.Range("F30").Value = 167
$ws.Range("F31").Value = 107
$ws.Range("F32").Value = 938
$ws.Range("F33").Value = 1289
$ws.Range("F34").Value = 102
$ws.Range("F35").Value = 138
$ws.Range("F36").Value = 1515
$ws.Range("F37").Value = 171
$ws.Range("F38").Value = 1556
$ws.Range("F39").Value = 113
$ws.Range("F40").Value = 998
$ws.Range("F41").Value = 1081
$ws.Range("F42").Value = 576
$ws.Range("F44").Value = 135
$ws.Range("F45").Value = 2641
$ws.Range("F47").Value = 223
$ws.Range("F48").Value = 390
$ws.Range("F49").Value = 3823

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 1117

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 3441

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 3441
$ws.Range("F3").Value = 37
$ws.Range("F4").Value = 821
$ws.Range("F5").Value = 2468
$ws.Range("F6").Value = 524
$ws.Range("F7").Value = 447
$ws.Range("F8").Value = 269
$ws.Range("F9").Value = 1117
$ws.Range("F12").Value = 1147
$ws.Range("F13").Value = 507
$ws.Range("F14").Value = 270
$ws.Range("F17").Value = 5250
$ws.Range("F18").Value = 1503
$ws.Range("F19").Value = 4286
$ws.Range("F20").Value = 5680
$ws.Range("F21").Value = 136
$ws.Range("F22").Value = 1001
$ws.Range("F23").Value = 615
$ws.Range("F24").Value = 3554
$ws.Range("F25").Value = 436
$ws.Range("F27").Value = 167
$ws.Range("F28").Value = 107
$ws.Range("F29").Value = 938
$ws.Range("F30").Value = 1289
$ws.Range("F31").Value = 102
$ws.Range("F32").Value = 138
$ws.Range("F33").Value = 1515
$ws.Range("F34").Value = 171
$ws.Range("F35").Value = 1556
$ws.Range("F37").Value = 998
$ws.Range("F39").Value = 576
$ws.Range("F43").Value = 2641
$ws.Range("F46").Value = 223
$ws.Range("F47").Value = 390
$ws.Range("F49").Value = 3823
